$wb = $excel.ActiveWorkbook

# Avoid interactive prompts when deleting a sheet.
$excel.DisplayAlerts = $false

# Rename "Paineis DARQ" -> "PAINEIS DARQ"
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"

# Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Remove the "Desarquivamentos Pendentes" sheet entirely.
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete() | Out-Null

# Deleting a sheet can shift which tab is active; restore the original
# active sheet ("PAINEIS DARQ", formerly "Paineis DARQ").
$wb.Worksheets.Item("PAINEIS DARQ").Activate()

$excel.DisplayAlerts = $true
